# Add two new columns to the sheet: I ("I0") and J ("IF"), matching the
# existing header style and populating the same number of data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold font, border, center/top alignment) from the
# existing "IP" header (H1) onto the two new header cells so they match
# the rest of the header row's style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data rows (2-32): column I = I0, column J = IF ---
$values = @(
    @(6, 6),
    @(6, 6),
    @(6, 6),
    @(8, 8),
    @(7, 8),
    @(6, 7),
    @(8, 8),
    @(9, 9),
    @(8, 8),
    @(1, 1),
    @(7, 7),
    @(6, 6),
    @(7, 7),
    @(8, 8),
    @(4, 5),
    @(5, 6),
    @(8, 8),
    @(6, 6),
    @(7, 7),
    @(6, 7),
    @(7, 7),
    @(6, 6),
    @(9, 9),
    @(6, 6),
    @(6, 6),
    @(8, 8),
    @(6, 6),
    @(5, 5),
    @(5, 5),
    @(5, 5),
    @(3, 3)
)

$row = 2
foreach ($pair in $values) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
